$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "books"
$ws.Name = "books"

# Append the new book catalogue rows (7-26).
# Values are written in the same cell-by-cell order used when the
# original data was authored, so the shared-string table comes out
# in the same sequence as the source workbook.
$ws.Cells.Item(7, 1).Value = 'The Sisters Brothers'
$ws.Cells.Item(7, 2).Value = 'Patrick DeWitt'
$ws.Cells.Item(8, 1).Value = 'The Essex Serpent'
$ws.Cells.Item(7, 3).Value = 'Western;Literature;Historical'
$ws.Cells.Item(8, 3).Value = 'Literature;Historical'
$ws.Cells.Item(9, 1).Value = 'Why I’m No Longer Talking to White People About Race'
$ws.Cells.Item(9, 2).Value = 'Reni Eddo-Lodge'
$ws.Cells.Item(9, 3).Value = 'Non-Fiction;Social Sciences'
$ws.Cells.Item(8, 2).Value = 'Sarah Perry'
$ws.Cells.Item(10, 1).Value = 'Cujo'
$ws.Cells.Item(10, 2).Value = 'Stephen King'
$ws.Cells.Item(10, 3).Value = 'Horror'
$ws.Cells.Item(11, 1).Value = 'Blackbirds'
$ws.Cells.Item(11, 2).Value = 'Chuck Wendig'
$ws.Cells.Item(11, 3).Value = 'Horror; Thriller'
$ws.Cells.Item(12, 1).Value = 'Hollow Things'
$ws.Cells.Item(12, 2).Value = 'T.S. King'
$ws.Cells.Item(12, 3).Value = 'Horror'
$ws.Cells.Item(13, 1).Value = 'Heart Shaped Box'
$ws.Cells.Item(13, 2).Value = 'Joe Hill'
$ws.Cells.Item(13, 3).Value = 'Horror'
$ws.Cells.Item(14, 1).Value = 'Buddha Da'
$ws.Cells.Item(14, 2).Value = 'Anne Donovan'
$ws.Cells.Item(15, 1).Value = 'Trainspotting'
$ws.Cells.Item(15, 2).Value = 'Irvine Welsh'
$ws.Cells.Item(16, 1).Value = 'The Crow Road'
$ws.Cells.Item(16, 2).Value = 'Iain Banks'
$ws.Cells.Item(14, 3).Value = 'Literature; Scottish'
$ws.Cells.Item(16, 3).Value = 'Literature; Scottish'
$ws.Cells.Item(15, 3).Value = 'Literature; Scottish'
$ws.Cells.Item(17, 1).Value = 'Klara and the Sun'
$ws.Cells.Item(17, 2).Value = 'Kazuo Ishiguro'
$ws.Cells.Item(17, 3).Value = 'Literature'
$ws.Cells.Item(18, 1).Value = 'One: Pot, Pan, Planet'
$ws.Cells.Item(18, 2).Value = 'Anna Jones'
$ws.Cells.Item(18, 3).Value = 'Non-Fiction;Cookery'
$ws.Cells.Item(19, 1).Value = 'The Midnight Library'
$ws.Cells.Item(19, 2).Value = 'Matt Haig'
$ws.Cells.Item(19, 3).Value = 'Modern Fiction'
$ws.Cells.Item(20, 1).Value = 'Acts of Desperation'
$ws.Cells.Item(20, 2).Value = 'Megan Nolan'
$ws.Cells.Item(20, 3).Value = 'Modern Fiction'
$ws.Cells.Item(21, 1).Value = 'Transcendant Kingdom'
$ws.Cells.Item(21, 2).Value = 'Yaa Gyasi'
$ws.Cells.Item(21, 3).Value = 'Modern Fiction'
$ws.Cells.Item(22, 1).Value = 'Difficult Women'
$ws.Cells.Item(22, 2).Value = 'Helen Lewis'
$ws.Cells.Item(22, 3).Value = 'Social Sciences'
$ws.Cells.Item(23, 1).Value = 'With These Hands'
$ws.Cells.Item(23, 2).Value = 'Pam Ayres'
$ws.Cells.Item(23, 3).Value = 'Biography'
$ws.Cells.Item(24, 1).Value = 'Empireland'
$ws.Cells.Item(24, 2).Value = 'Sathnam Sanghera'
$ws.Cells.Item(24, 3).Value = 'History; British History'
$ws.Cells.Item(25, 1).Value = 'The Thursday Murder Club'
$ws.Cells.Item(25, 2).Value = 'Richard Osman'
$ws.Cells.Item(25, 3).Value = 'Modern Fiction; Crime; Thriller'
$ws.Cells.Item(26, 1).Value = 'Luster'
$ws.Cells.Item(26, 2).Value = 'Raven Leilani'
$ws.Cells.Item(26, 3).Value = 'Modern Fiction'

# Fill in the Price column for the new rows.
for ($r = 7; $r -le 26; $r++) {
    $ws.Cells.Item($r, 4).Value = 8.99
}

# Move the active selection to match the post-edit state.
$ws.Range("E34").Select()

